$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatted rows 17-19 by copying format+value from row 16, then overwrite with correct content
$ws.Range("A16:B16").Copy($ws.Range("A17"))
$ws.Range("A16:B16").Copy($ws.Range("A18"))
$ws.Range("A16:B16").Copy($ws.Range("A19"))

# Set A (index) and B (scheme name) labels for rows 10-19
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"

# Set numeric intensity values C:M for rows 10-19
$ws.Range("C10").Value = 0.9842483788308727
$ws.Range("D10").Value = 1.021358771308097
$ws.Range("E10").Value = 0.9860784505992501
$ws.Range("F10").Value = 0.9842483788308727
$ws.Range("G10").Value = 1.01029065399248
$ws.Range("H10").Value = 0.9735217150096211
$ws.Range("I10").Value = 0.9831372741286619
$ws.Range("J10").Value = 1.021358771308097
$ws.Range("K10").Value = 1.003718610953674
$ws.Range("L10").Value = 0.993983494892273
$ws.Range("M10").Value = 0.9931058739781639
$ws.Range("C11").Value = 0.9992114341895088
$ws.Range("D11").Value = 0.9718420115976726
$ws.Range("E11").Value = 1.000089702461063
$ws.Range("F11").Value = 0.9992114341895088
$ws.Range("G11").Value = 0.980062931742587
$ws.Range("H11").Value = 1.013465366313336
$ws.Range("I11").Value = 0.9989068172540068
$ws.Range("J11").Value = 0.9718420115976726
$ws.Range("K11").Value = 0.9859658570293679
$ws.Range("L11").Value = 0.9925886456094383
$ws.Range("M11").Value = 0.9939297105930289
$ws.Range("C12").Value = 0.9991983676213081
$ws.Range("D12").Value = 0.9719958760791825
$ws.Range("E12").Value = 1.000051165419836
$ws.Range("F12").Value = 0.9991983676213081
$ws.Range("G12").Value = 0.9801654765517988
$ws.Range("H12").Value = 1.013398184189675
$ws.Range("I12").Value = 0.9988718232249534
$ws.Range("J12").Value = 0.9719958760791825
$ws.Range("K12").Value = 0.9860235207495091
$ws.Range("L12").Value = 0.9926109441854085
$ws.Range("M12").Value = 0.9939468155144588
$ws.Range("C13").Value = 0.9992235183805473
$ws.Range("D13").Value = 0.9719135306014605
$ws.Range("E13").Value = 1.00006959388233
$ws.Range("F13").Value = 0.9992235183805473
$ws.Range("G13").Value = 0.9801048495546446
$ws.Range("H13").Value = 1.013481658129584
$ws.Range("I13").Value = 0.9988940463938802
$ws.Range("J13").Value = 0.9719135306014605
$ws.Range("K13").Value = 0.9859915622418953
$ws.Range("L13").Value = 0.9926075403112213
$ws.Range("M13").Value = 0.9939478661570744
$ws.Range("C14").Value = 0.9750000000000002
$ws.Range("D14").Value = 1.059372
$ws.Range("E14").Value = 0.9745839999999993
$ws.Range("F14").Value = 0.9750000000000002
$ws.Range("G14").Value = 1.031044
$ws.Range("H14").Value = 0.9383199999999999
$ws.Range("I14").Value = 0.9745840000000002
$ws.Range("J14").Value = 1.059372
$ws.Range("K14").Value = 1.016978
$ws.Range("L14").Value = 0.9959889999999999
$ws.Range("M14").Value = 0.9921506666666665
$ws.Range("C15").Value = 0.96
$ws.Range("D15").Value = 1.11
$ws.Range("E15").Value = 0.96
$ws.Range("F15").Value = 0.96
$ws.Range("G15").Value = 1.061250000000001
$ws.Range("H15").Value = 0.89
$ws.Range("I15").Value = 0.96
$ws.Range("J15").Value = 1.11
$ws.Range("K15").Value = 1.035
$ws.Range("L15").Value = 0.9975000000000001
$ws.Range("M15").Value = 0.9902083333333335
$ws.Range("C16").Value = 0.9748132051968023
$ws.Range("D16").Value = 1.062359936819193
$ws.Range("E16").Value = 0.9744879730687999
$ws.Range("F16").Value = 0.9748132051968023
$ws.Range("G16").Value = 1.033198513663996
$ws.Range("H16").Value = 0.9349155694592045
$ws.Range("I16").Value = 0.9747430590464036
$ws.Range("J16").Value = 1.062359936819193
$ws.Range("K16").Value = 1.018423954943997
$ws.Range("L16").Value = 0.9966185800703994
$ws.Range("M16").Value = 0.9924197095423999
$ws.Range("C17").Value = 0.99479708822698
$ws.Range("D17").Value = 0.9945805313864726
$ws.Range("E17").Value = 0.9944214469720416
$ws.Range("F17").Value = 0.99479708822698
$ws.Range("G17").Value = 0.9949477028747656
$ws.Range("H17").Value = 0.9942016177834837
$ws.Range("I17").Value = 0.9946708432414045
$ws.Range("J17").Value = 0.9945805313864726
$ws.Range("K17").Value = 0.9945009891792571
$ws.Range("L17").Value = 0.9946490387031186
$ws.Range("M17").Value = 0.9946032050808581
$ws.Range("C18").Value = 0.9956095208745412
$ws.Range("D18").Value = 0.9895800295334213
$ws.Range("E18").Value = 0.9946006024546515
$ws.Range("F18").Value = 0.9956095208745412
$ws.Range("G18").Value = 0.9915404380170352
$ws.Range("H18").Value = 0.9960858699776257
$ws.Range("I18").Value = 0.9947088487688238
$ws.Range("J18").Value = 0.9895800295334213
$ws.Range("K18").Value = 0.9920903159940364
$ws.Range("L18").Value = 0.9938499184342888
$ws.Range("M18").Value = 0.9936875516043499
$ws.Range("C19").Value = 0.9978601197461112
$ws.Range("D19").Value = 0.9814681602980121
$ws.Range("E19").Value = 0.9975773319353027
$ws.Range("F19").Value = 0.9978601197461112
$ws.Range("G19").Value = 0.9868877929607551
$ws.Range("H19").Value = 1.003620805669126
$ws.Range("I19").Value = 0.9979561835851133
$ws.Range("J19").Value = 0.9814681602980121
$ws.Range("K19").Value = 0.9895227461166574
$ws.Range("L19").Value = 0.9936914329313842
$ws.Range("M19").Value = 0.9942283990324036
